# Applies the cryptos.xlsx crypto price/volume data refresh
# described by the commit "Updated cryptos list on Mon Jul 24 22:38:53 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force these Price cells to stay text (otherwise Excel would parse the new
# numeric-looking value as a real number instead of keeping it as a string)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.182.99"
$ws.Range("E2").Value = "  -2.81%  "

$ws.Range("D3").Value = "1.847.56"

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "0.7032"
$ws.Range("E5").Value = "  -4.69%  "

$ws.Range("D6").Value = "238.68"
$ws.Range("E6").Value = "  -1.37%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "0.3046"
$ws.Range("E8").Value = "  -3.61%  "

$ws.Range("D9").Value = "0.07407"
$ws.Range("E9").Value = "  +3.27%  "

$ws.Range("E10").Value = "  -5.04%  "

$ws.Range("D11").Value = "0.08133"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "0.7273"
$ws.Range("E12").Value = "  -3.79%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.848.82"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").Value = "5.217"
$ws.Range("E14").Value = "  -3.43%  "

$ws.Range("D15").Value = "88.79"
$ws.Range("E15").Value = "  -4.11%  "

$ws.Range("D16").Value = "29.201.77"
$ws.Range("E16").Value = "  -2.87%  "

$ws.Range("D17").Value = "5.762"
$ws.Range("E17").Value = "  -6.23%  "

$ws.Range("D18").Value = "238.54"
$ws.Range("E18").Value = "  -4.45%  "

$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").Value = "0.000007647"
$ws.Range("E20").Value = "  -2.57%  "

$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.099.94"
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").Value = "7.603"
$ws.Range("E24").Value = "  -3.87%  "

$ws.Range("D25").Value = "8.996"
$ws.Range("E25").Value = "  -2.90%  "

$ws.Range("D26").Value = "160.51"
$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("D27").Value = "0.1450"
$ws.Range("E27").Value = "  -7.78%  "

$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  -3.15%  "

$ws.Range("D29").Value = "1.969"
$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("D30").Value = "1.397"
$ws.Range("E30").Value = "  -5.09%  "

$ws.Range("D31").Value = "4.518"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").Value = "1.490"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").Value = "3.990"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("D34").Value = "0.05188"
$ws.Range("E34").Value = "  -2.55%  "

$ws.Range("D35").Value = "1.187"
$ws.Range("E35").Value = "  -4.93%  "

$ws.Range("E36").Value = "  +2.94%  "

$ws.Range("D37").Value = "0.7029"
$ws.Range("E37").Value = "  -8.37%  "

$ws.Range("D38").Value = "2.666"
$ws.Range("E38").Value = "  -2.23%  "

$ws.Range("D39").Value = "0.01871"
$ws.Range("E39").Value = "  -4.52%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.676"
$ws.Range("E40").Value = "  -2.94%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9606"
$ws.Range("E41").Value = "  +8.96%  "

$ws.Range("D42").Value = "6.004"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("D43").Value = "0.4296"
$ws.Range("E43").Value = "  -5.74%  "

$ws.Range("D44").Value = "1.072.00"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("D45").Value = "70.30"
$ws.Range("E45").Value = "  -2.84%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").Value = "102.83"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("E48").Value = "  -6.09%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.990.86"
$ws.Range("E49").Value = "  -3.58%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.037"
$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("D51").Value = "9.115"
$ws.Range("E51").Value = "  -4.39%  "
